$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts D:K -> E:L), matching the new reporting period added to each table
$ws.Columns("D:D").Insert()

# Copy column E formatting (the old column D, now shifted) onto new column D so styles/number formats match
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(-4122)

# Populate the new column D with the new period values
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 58400
$ws.Range("D9").Value2 = 24200
$ws.Range("D10").Value2 = 34200
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("D17").Value2 = 26500
$ws.Range("D18").Value2 = 31900
$ws.Range("D20").Value2 = -26100
$ws.Range("D21").Value2 = "NA"
$ws.Range("D22").Value2 = 0
$ws.Range("D23").Value2 = 5800
$ws.Range("D24").Value2 = 0
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 5800
$ws.Range("D27").Value2 = 5800
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = 26100
$ws.Range("D33").Value2 = 5800
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 5800
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 3700
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 7800
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 0
$ws.Range("D46").Value2 = 0
$ws.Range("D47").Value2 = 553600
$ws.Range("D48").Value2 = 0
$ws.Range("D49").Value2 = 0
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 14000
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 579800
$ws.Range("D57").Value2 = 2400
$ws.Range("D58").Value2 = "NA"
$ws.Range("D59").Value2 = 4900
$ws.Range("D60").Value2 = 0
$ws.Range("D61").Value2 = 320000
$ws.Range("D62").Value2 = 0
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 321100
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = "NA"
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 258800
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 5800
$ws.Range("D83").Value2 = "NA"
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = -55100
$ws.Range("D91").Value2 = 0
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = "NA"
$ws.Range("D96").Value2 = -28200
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = 65700
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = 10500

# A few historical cells were revised to "NA" as part of this data refresh (beyond the simple column shift)
$ws.Range("E58:J58").Value2 = "NA"
$ws.Range("J83").Value2 = "NA"
$ws.Range("E101:J101").Value2 = "NA"
